# Add a duplicate "Actual" column (D) next to the existing Actual/Nodes
# columns on the second worksheet, mirroring column B's values into D.
# This reproduces the "bad generator" edge-label-partition column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header: D3 should reuse the existing "Actual" shared string (same as B3).
$ws.Cells.Item(3, 4).Value = $ws.Cells.Item(3, 2).Value2

# Data rows: D4:D23 mirror B4:B23 (values + number format/style).
for ($r = 4; $r -le 23; $r++) {
    $src = $ws.Cells.Item($r, 2)   # column B
    $dst = $ws.Cells.Item($r, 4)   # column D
    $dst.Value = $src.Value2
    $dst.NumberFormat = $src.NumberFormat
}

# Match the recorded selection left behind after the edit.
[void]$ws.Range("C3:D23").Select()

# Force a portrait page setup (as recorded after the edit).
$ws.PageSetup.Orientation = 1
